$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '33.787.06'
$ws.Range("E2").Value = '  +9.97%  '
$ws.Range("D3").Value = '1.803.64'
$ws.Range("E3").Value = '  +6.99%  '
$ws.Range("E4").Value = '  +0.08%  '
$ws.Range("D5").Value = '227.21'
$ws.Range("E5").Value = '  +3.07%  '
$ws.Range("D6").Value = '0.537'
$ws.Range("E6").Value = '  +1.91%  '
$ws.Range("E7").Value = '  +0.10%  '
$ws.Range("D8").Value = '30.89'
$ws.Range("E8").Value = '  +1.34%  '
$ws.Range("D9").Value = '47.24'
$ws.Range("E9").Value = '  +6.44%  '
$ws.Range("E10").Value = '  +5.20%  '
$ws.Range("E11").Value = '  +5.43%  '
$ws.Range("D12").Value = '0.0928'
$ws.Range("E12").Value = '  +2.19%  '
$ws.Range("D13").Value = '2.064.94'
$ws.Range("E13").Value = '  +7.06%  '
$ws.Range("D14").Value = '1.810.78'
$ws.Range("E14").Value = '  +7.50%  '
$ws.Range("D15").Value = '0.632'
$ws.Range("E15").Value = '  +1.74%  '
$ws.Range("D16").Value = '33.730.57'
$ws.Range("E16").Value = '  +9.76%  '
$ws.Range("E17").Value = '  -3.67%  '
$ws.Range("E18").Value = '  +5.82%  '
$ws.Range("E19").Value = '  +3.79%  '
$ws.Range("D20").Value = '254.35'
$ws.Range("E20").Value = '  +3.27%  '
$ws.Range("E21").Value = '  +3.17%  '
$ws.Range("D22").Value = '1.00'
$ws.Range("E22").Value = '  +0.08%  '
$ws.Range("E23").Value = '  +1.75%  '
$ws.Range("E24").Value = '  -0.07%  '
$ws.Range("E25").Value = '  +0.80%  '
$ws.Range("D26").Value = '157.67'
$ws.Range("E26").Value = '  -0.32%  '
$ws.Range("D27").Value = '16.35'
$ws.Range("E27").Value = '  +2.95%  '
$ws.Range("E28").Value = '  +2.33%  '
$ws.Range("E29").Value = '  +4.56%  '
$ws.Range("B30").Value = 'MinaProtocolToken'
$ws.Range("C30").Value = 'https://coinranking.com/coin/J7st_qGwz+minaprotocoltoken-mina'
$ws.Range("D30").Value = '2.11'
$ws.Range("E30").Value = '  +417.22%  '
$ws.Range("B31").Value = 'BinanceUSD'
$ws.Range("C31").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D31").Value = '1.00'
$ws.Range("E31").Value = '  +0.11%  '
$ws.Range("D32").Value = '3.79'
$ws.Range("E32").Value = '  +8.98%  '
$ws.Range("B33").Value = 'Hedera'
$ws.Range("C33").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D33").Value = '0.0507'
$ws.Range("E33").Value = '  +1.65%  '
$ws.Range("B34").Value = 'PancakeSwap'
$ws.Range("C34").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D34").Value = '1.19'
$ws.Range("E34").Value = '  +4.83%  '
$ws.Range("D35").Value = '3.48'
$ws.Range("E35").Value = '  +5.56%  '
$ws.Range("D36").Value = '1.529.22'
$ws.Range("E36").Value = '  +1.03%  '
$ws.Range("E37").Value = '  +1.97%  '
$ws.Range("D38").Value = '1.06'
$ws.Range("E38").Value = '  +2.98%  '
$ws.Range("B39").Value = 'Aave'
$ws.Range("C39").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D39").Value = '83.18'
$ws.Range("E39").Value = '  -1.57%  '
$ws.Range("B40").Value = 'VeChain'
$ws.Range("C40").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D40").Value = '0.0184'
$ws.Range("E40").Value = '  +3.42%  '
$ws.Range("D41").Value = '0.612'
$ws.Range("E41").Value = '  +4.12%  '
$ws.Range("D42").Value = '2.79'
$ws.Range("E42").Value = '  +2.39%  '
$ws.Range("E43").Value = '  +1.13%  '
$ws.Range("E44").Value = '  +7.14%  '
$ws.Range("E45").Value = '  +5.48%  '
$ws.Range("E46").Value = '  +3.91%  '
$ws.Range("E47").Value = '  +3.85%  '
$ws.Range("D48").Value = '1.952.99'
$ws.Range("E48").Value = '  +7.01%  '
$ws.Range("E49").Value = '  +0.03%  '
$ws.Range("E50").Value = '  +3.63%  '
$ws.Range("D51").Value = '51.77'
$ws.Range("E51").Value = '  -0.63%  '
